$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data set added three new weekly observations for
# "Agrícola del Norte S.A. de Arica - Choclo" dated 2022-07-05 (serial 44747),
# inserted right after the existing row for serial 44707 (row 676) and
# before the existing run of rows that started at serial 44421 (old row 677).
# Insert 3 blank rows at 677:679 - this shifts the old rows 677-691 down to 680-694.
$ws.Rows("677:679").Insert()

# New row 677: Lluteño / Primera
$ws.Cells.Item(677, 1).Value = 1
$ws.Cells.Item(677, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(677, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(677, 4).Value = 44747
$ws.Cells.Item(677, 5).Value = 15
$ws.Cells.Item(677, 6).Value = 100112024
$ws.Cells.Item(677, 7).Value = "Choclo"
$ws.Cells.Item(677, 8).Value = "Lluteño"
$ws.Cells.Item(677, 9).Value = "Primera"
$ws.Cells.Item(677, 10).Value = 50
$ws.Cells.Item(677, 11).Value = 35000
$ws.Cells.Item(677, 12).Value = 36000
$ws.Cells.Item(677, 13).Value = 35500
$ws.Cells.Item(677, 14).Value = "$/saco 50 unidades"
$ws.Cells.Item(677, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(677, 16).Value = 710
$ws.Cells.Item(677, 17).Value = 50
$ws.Cells.Item(677, 18).Value = "Hortaliza"

# New row 678: Lluteño / Segunda
$ws.Cells.Item(678, 1).Value = 1
$ws.Cells.Item(678, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(678, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(678, 4).Value = 44747
$ws.Cells.Item(678, 5).Value = 15
$ws.Cells.Item(678, 6).Value = 100112024
$ws.Cells.Item(678, 7).Value = "Choclo"
$ws.Cells.Item(678, 8).Value = "Lluteño"
$ws.Cells.Item(678, 9).Value = "Segunda"
$ws.Cells.Item(678, 10).Value = 50
$ws.Cells.Item(678, 11).Value = 30000
$ws.Cells.Item(678, 12).Value = 32000
$ws.Cells.Item(678, 13).Value = 31000
$ws.Cells.Item(678, 14).Value = "$/saco 75 unidades"
$ws.Cells.Item(678, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(678, 16).Value = 413
$ws.Cells.Item(678, 17).Value = 75
$ws.Cells.Item(678, 18).Value = "Hortaliza"

# New row 679: Lluteño / Tercera
$ws.Cells.Item(679, 1).Value = 1
$ws.Cells.Item(679, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(679, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(679, 4).Value = 44747
$ws.Cells.Item(679, 5).Value = 15
$ws.Cells.Item(679, 6).Value = 100112024
$ws.Cells.Item(679, 7).Value = "Choclo"
$ws.Cells.Item(679, 8).Value = "Lluteño"
$ws.Cells.Item(679, 9).Value = "Tercera"
$ws.Cells.Item(679, 10).Value = 40
$ws.Cells.Item(679, 11).Value = 25000
$ws.Cells.Item(679, 12).Value = 26000
$ws.Cells.Item(679, 13).Value = 25500
$ws.Cells.Item(679, 14).Value = "$/saco 100 unidades"
$ws.Cells.Item(679, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(679, 16).Value = 255
$ws.Cells.Item(679, 17).Value = 100
$ws.Cells.Item(679, 18).Value = "Hortaliza"
